# Auto-committed on 2022/08/15 week-Mon 18:59:13.58
# Adds a new "OtherNote" (其他說明) field row to the BatxDetail (DBD) field
# list, right after "ProcNote", and appends a new code ("99:暫收沖正") to the
# RepayCode note text in row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# --- 1. Append the new BatchRepayCode entry to the row-13 note (G13). ---
$g13 = $ws.Range("G13").Value()
$ws.Range("G13").Value = $g13 + "`n99:暫收沖正"
$ws.Rows.Item(13).RowHeight = 194.4

# --- 2. Insert a new row for the "OtherNote" field, right after "ProcNote"
#        (row 28), pushing rows 29:37 down to 30:38. ---
$ws.Rows.Item(29).Insert()

# Copy the formatting (borders/fonts/number-format/alignment) from the
# "ProcNote" row above so the new row matches its neighbours exactly.
$ws.Range("A28:G28").Copy()
$ws.Range("A29:G29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Populate the new row's values. ---
$ws.Range("A29").Value = 20
$ws.Range("B29").Value = "OtherNote"
$ws.Range("C29").Value = "其他說明"
$ws.Range("D29").Value = "NVARCHAR2"
$ws.Range("E29").Value = 2000
$ws.Range("G29").Value = "jsonformat"

# Renumber the SEQ column (A) for every row pushed down by the insert, since
# those were literal numbers, not a formula.
for ($r = 30; $r -le 38; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 9
}

# --- 4. Restore the view: scroll/select around the newly-inserted row. ---
$ws.Range("G28").Select()

Write-Host "BatxDetail: inserted OtherNote row; updated RepayCode note."
